$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching H1 formatting (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I2:J80 with the new data values
$ijData = @{
  2 = @(9, 10)
  3 = @(9, 9)
  4 = @(6, 6)
  5 = @(8, 8)
  6 = @(9, 9)
  7 = @(4, 4)
  8 = @(7, 7)
  9 = @(5, 6)
  10 = @(6, 6)
  11 = @(7, 7)
  12 = @(6, 6)
  13 = @(8, 8)
  14 = @(8, 8)
  15 = @(9, 9)
  16 = @(9, 9)
  17 = @(1, 1)
  18 = @(8, 9)
  19 = @(6, 7)
  20 = @(8, 8)
  21 = @(7, 7)
  22 = @(8, 8)
  23 = @(8, 8)
  24 = @(7, 7)
  25 = @(6, 7)
  26 = @(7, 7)
  27 = @(8, 8)
  28 = @(9, 9)
  29 = @(6, 6)
  30 = @(6, 6)
  31 = @(7, 7)
  32 = @(6, 6)
  33 = @(7, 7)
  34 = @(8, 8)
  35 = @(6, 6)
  36 = @(7, 7)
  37 = @(5, 6)
  38 = @(6, 6)
  39 = @(8, 8)
  40 = @(5, 5)
  41 = @(1, 2)
  42 = @(8, 8)
  43 = @(7, 7)
  44 = @(9, 9)
  45 = @(6, 6)
  46 = @(8, 8)
  47 = @(7, 7)
  48 = @(7, 7)
  49 = @(7, 7)
  50 = @(7, 7)
  51 = @(8, 8)
  52 = @(5, 6)
  53 = @(6, 6)
  54 = @(6, 6)
  55 = @(6, 6)
  56 = @(3, 4)
  57 = @(6, 6)
  58 = @(6, 6)
  59 = @(6, 6)
  60 = @(4, 5)
  61 = @(6, 6)
  62 = @(8, 8)
  63 = @(8, 8)
  64 = @(7, 7)
  65 = @(8, 8)
  66 = @(5, 5)
  67 = @(6, 6)
  68 = @(7, 7)
  69 = @(6, 7)
  70 = @(7, 8)
  71 = @(7, 7)
  72 = @(8, 8)
  73 = @(8, 8)
  74 = @(8, 8)
  75 = @(7, 7)
  76 = @(8, 8)
  77 = @(6, 6)
  78 = @(6, 6)
  79 = @(7, 7)
  80 = @(4, 4)
}

foreach ($r in $ijData.Keys) {
  $vals = $ijData[$r]
  $ws.Cells.Item($r, 9).Value = $vals[0]
  $ws.Cells.Item($r, 10).Value = $vals[1]
}
